# Replace the four "letter of transmittal" bullet paragraphs with the
# updated set of paragraphs (bullets kept, each followed by a new prose
# paragraph with the actual report content; trailing "An ..." bullet is
# merged into a single run and followed by one more new paragraph).
$d = $word.ActiveDocument

$newParagraphsXml = @'
    <w:p w:rsidR="00F00151" w:rsidRDefault="00F00151" w:rsidP="00F00151">
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hAnsi="Symbol"/>
        </w:rPr>
        <w:t></w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t>A</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t xml:space="preserve"> brief explanation of why it was written (Who assigned it? When? Why? How?)</w:t>
      </w:r>
    </w:p>
    <w:p w:rsidR="00F00151" w:rsidRDefault="00F00151" w:rsidP="00F00151">
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t xml:space="preserve">This report was written to provide context for the evaluation of infrared cameras for the detection of obstacles at sea, and to indicate the next steps for UBC Sailbot team in the goal of avoiding hazards during their upcoming autonomous crossing of the Atlantic </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t>Ocean</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t xml:space="preserve"> The Sailbot team has been planning to attempt this crossing since </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t>Fall</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t xml:space="preserve"> 2014, and our group has assisted with this project from January to the end of March 2015.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hAnsi="Symbol"/>
        </w:rPr>
        <w:t></w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t>A</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t xml:space="preserve"> one- to -three-sentence synopsis of the main theme of the report</w:t>
      </w:r>
    </w:p>
    <w:p w:rsidR="00F00151" w:rsidRDefault="00F00151" w:rsidP="00F00151">
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t>Contained in the report is an outline of the equipment developed for the protection during testing of a small infrared sensor</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t>, and consideration of</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t xml:space="preserve"> the effectiveness of a FLIR Lepton sensor for the detection of floating obstacles on the open ocean.</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hAnsi="Symbol"/>
        </w:rPr>
        <w:t></w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t>A</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t xml:space="preserve"> brief note on any particular features or sections of the report that may be of special interest to the recipient</w:t>
      </w:r>
    </w:p>
    <w:p w:rsidR="00F66B4C" w:rsidRPr="00F00151" w:rsidRDefault="00F00151" w:rsidP="00F00151">
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hAnsi="Symbol"/>
        </w:rPr>
        <w:t></w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t>An indication of what the recipient is to do with the report, particularly if the report is not directed primarily to him or her.  In other words, should the recipient respond in some way to it, or simply be aware and keep it on file?</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="NormalWeb"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
        </w:rPr>
        <w:t>The recipient of this report may keep it on file.</w:t>
      </w:r>
    </w:p>
'@

$pkgPrefix = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgSuffix = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xml = $pkgPrefix + $newParagraphsXml + $pkgSuffix

# The four bullet paragraphs to be replaced are Paragraphs 6-9.
$firstBullet = $d.Paragraphs(6)
$lastBullet = $d.Paragraphs(9)
$block = $d.Range($firstBullet.Range.Start, $lastBullet.Range.End)
$insertAt = $block.Start
$block.Delete()

$target = $d.Range($insertAt, $insertAt)
$target.InsertXML($xml)

Write-Output "Letter of transmittal bullets replaced."
